$d = $word.ActiveDocument

# Appends a new paragraph right after the paragraph owning $afterRange, sets
# its full text in one shot (so it inherits the surrounding run formatting
# cleanly), and then splits it into multiple <w:r> runs at the given
# 0-based character offsets (measured from the start of the paragraph's
# text) without losing formatting on the later pieces.
function Add-Para($afterRange, [string]$text, [int[]]$splitOffsets) {
    $afterRange.Collapse(0)
    $afterRange.InsertParagraphAfter()
    $count = $d.Paragraphs.Count
    $p = $d.Paragraphs.Item($count)
    $r = $p.Range
    if ($text -ne $null -and $text -ne "") {
        $r.Text = $text
    }
    if ($splitOffsets -ne $null -and $splitOffsets.Count -gt 0) {
        $pstart = $p.Range.Start
        $pend = $p.Range.End
        $sorted = $splitOffsets | Sort-Object -Descending
        foreach ($off in $sorted) {
            $b = $d.Range($pstart + $off, $pend - 1)
            $b.Font.Bold = 1
            $b.Font.Bold = 0
        }
    }
    return $d.Paragraphs.Item($d.Paragraphs.Count).Range
}

# 1) Rewrite the last paragraph's text ("---" -> "Câu 4: ...")
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$r1 = $lastP.Range
$r1.Text = "Câu 4: Phần hướng dẫn chạy file "

# 2) "Biên dịch và chạy (C++)"
$r2 = Add-Para $r1 "Biên dịch và chạy (C++)" $null

# 3) "- Biên dịch:"
$r3 = Add-Para $r2 "- Biên dịch:" $null

# 4) multi-run paragraph: "  g++ -std=c++17 " + "cau2" + ".cpp -o " + "cau2"
$seg4a = "  g++ -std=c++17 "
$seg4b = "cau2"
$seg4c = ".cpp -o "
$seg4d = "cau2"
$text4 = $seg4a + $seg4b + $seg4c + $seg4d
$off4b = $seg4a.Length
$off4c = $off4b + $seg4b.Length
$off4d = $off4c + $seg4c.Length
$r4 = Add-Para $r3 $text4 @($off4b, $off4c, $off4d)

# 5) "- Chạy:"
$r5 = Add-Para $r4 "- Chạy:" $null

# 6) "  ./" + "cau2.cpp"
$seg6a = "  ./"
$seg6b = "cau2.cpp"
$text6 = $seg6a + $seg6b
$off6b = $seg6a.Length
$r6 = Add-Para $r5 $text6 @($off6b)

# 7) empty paragraph
$r7 = Add-Para $r6 $null $null

# 8) "Chạy Shell script"
$r8 = Add-Para $r7 "Chạy Shell script" $null

# 9) "- Đảm bảo script có quyền thực thi:"
$r9 = Add-Para $r8 "- Đảm bảo script có quyền thực thi:" $null

# 10) "  chmod +x " + "cau2" + ".sh"
$seg10a = "  chmod +x "
$seg10b = "cau2"
$seg10c = ".sh"
$text10 = $seg10a + $seg10b + $seg10c
$off10b = $seg10a.Length
$off10c = $off10b + $seg10b.Length
$r10 = Add-Para $r9 $text10 @($off10b, $off10c)

# 11) "- Chạy:"
$r11 = Add-Para $r10 "- Chạy:" $null

# 12) "  ./" + "cau2" + ".sh"
$seg12a = "  ./"
$seg12b = "cau2"
$seg12c = ".sh"
$text12 = $seg12a + $seg12b + $seg12c
$off12b = $seg12a.Length
$off12c = $off12b + $seg12b.Length
$r12 = Add-Para $r11 $text12 @($off12b, $off12c)

Write-Output "Done. Paragraph count:"
Write-Output $d.Paragraphs.Count
